$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update market cap values (row order unchanged for most rows)
$ws.Range("C2").Value = 707764416716.2303
$ws.Range("C3").Value = 241861615405.5887
$ws.Range("C4").Value = 87347430651.04588
$ws.Range("C5").Value = 37376344664.52251
$ws.Range("C6").Value = 34322046171.68675

# Rows 7 and 8 swap: Solana moves up to row 7, USDC moves down to row 8
$ws.Range("A7").Value = "Solana"
$ws.Range("B7").Value = "SOL-USD"
$ws.Range("C7").Value = 26604916001.90222

$ws.Range("A8").Value = "USDC"
$ws.Range("B8").Value = "USDC-USD"
$ws.Range("C8").Value = 23931425825.1214

$ws.Range("C9").Value = 13037134434.13486
$ws.Range("C10").Value = 10583343498.86235
$ws.Range("C11").Value = 9261916358.333483
$ws.Range("C12").Value = 8763528480.611448
$ws.Range("C13").Value = 8318387938.914449
$ws.Range("C14").Value = 7889662081.524138

# Rows 15 and 16 swap: Avalanche moves up to row 15, Polkadot moves down to row 16
$ws.Range("A15").Value = "Avalanche"
$ws.Range("B15").Value = "AVAX-USD"
$ws.Range("C15").Value = 6959721541.2058

$ws.Range("A16").Value = "Polkadot"
$ws.Range("B16").Value = "DOT-USD"
$ws.Range("C16").Value = 6774478471.431178
